$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.446.48"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "3.775.23"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "596.56"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "169.11"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "3.773.87"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "6.54"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "36.87"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "4.411.04"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "3.771.08"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "68.462.20"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "18.19"
$ws.Range("E18").Value = "  -3.37%  "
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "10.94"
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("D22").Value = "468.51"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "0.704"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").Value = "85.07"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "0.0000144"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "10.19"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "3.922.75"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D34").Value = "30.11"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("D37").Value = "3.729.87"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "3.51"
$ws.Range("E39").Value = "  -9.42%  "
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D46").Value = "1.97"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "42.21"
$ws.Range("E48").Value = "  +8.52%  "
$ws.Range("D49").Value = "403.23"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "45.73"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "145.64"
$ws.Range("E51").Value = "  +3.06%  "
